# Insert a new data row at row 7 (pushing the existing rows 7-13 down to 8-14)
# and populate it with the new "Primera" quality entry (week of 2023-06-02).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(7).Insert()

$ws.Range("A7").Value = 10
$ws.Range("B7").Value = "Vega Modelo de Temuco"
$ws.Range("C7").Value = "La Araucanía"
$ws.Range("D7").Value = 45079
$ws.Range("E7").Value = 9
$ws.Range("F7").Value = "Fruta"
$ws.Range("G7").Value = 100104
$ws.Range("H7").Value = "Frutos de pepita"
$ws.Range("I7").Value = 100104005
$ws.Range("J7").Value = "Pera asiática"
$ws.Range("K7").Value = "Hosui"
$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 100
$ws.Range("N7").Value = 18000
$ws.Range("O7").Value = 18000
$ws.Range("P7").Value = 18000
$ws.Range("Q7").Value = "$/caja 18 kilos granel"
$ws.Range("R7").Value = "Región de O'Higgins"
$ws.Range("S7").Value = 1000
$ws.Range("T7").Value = 18
